$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy the header formatting
# (bold, bordered, centered style) from the existing H1 header cell.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J for rows 2-4 (plain numeric cells, no style)
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 4
